# feat: add 2022-Q1 data
#
# Before:  sheets = 2021-Q3, 2021-Q4, 总计
# After:   sheets = 2021-Q3, 2021-Q4, 2022-Q1, 总计
#
# The previously-last sheet ("总计") becomes "2022-Q1" and is populated with
# the per-fund holdings table for the new quarter. A brand new "总计" sheet
# is appended after it, holding the running summary-by-quarter table
# (now including the 2022-Q1 row).

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122
$xlCenter = -4108
$xlTop = -4160

# ---------------------------------------------------------------------
# 1. Duplicate the existing "总计" sheet so the new "总计" sheet inherits
#    identical page setup / sheet formatting boilerplate, then re-purpose
#    the original as "2022-Q1" and the copy as the new "总计".
# ---------------------------------------------------------------------
$quarterSheet = $wb.Worksheets.Item("总计")
$quarterSheet.Copy([System.Reflection.Missing]::Value, $quarterSheet)
$totalSheet = $wb.Worksheets.Item($quarterSheet.Index + 1)

$quarterSheet.Name = "2022-Q1"
$totalSheet.Name = "总计"

# ---------------------------------------------------------------------
# 2. Rebuild the "2022-Q1" sheet with the fund holdings table.
# ---------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$headerCols = @("B", "C", "D", "E", "F", "G", "H")

# Make sure every header cell (including the newly-added E:H columns)
# carries the same bold/centered/bordered style already used by B1:D1.
$quarterSheet.Range("B1").Copy()
$quarterSheet.Range("E1:H1").PasteSpecial($xlPasteFormats)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $quarterSheet.Range($headerCols[$i] + "1").Value = $headers[$i]
}

$fundRows = @(
    @("005585", "银河文体娱乐主题灵活配置混合", "5.54", "74.07", "5.54", "0.3069", 4),
    @("005381", "泰康睿利量化多策略混合A", "0.99", "93.49", "1.95", "0.0193", 8),
    @("005844", "东方人工智能主题混合", "0.40", "94.80", "3.70", "0.0148", 9),
    @("005382", "泰康睿利量化多策略混合C", "0.49", "93.49", "1.95", "0.0096", 8),
    @("003366", "浙商汇金中证转型成长指数", "0.09", "93.88", "1.24", "0.0011", 5)
)

# Column A keeps the bordered/centered style already present on A2; copy
# its formatting down onto the newly-needed A3:A6 cells.
$quarterSheet.Range("A2").Copy()
$quarterSheet.Range("A3:A6").PasteSpecial($xlPasteFormats)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $row = $i + 2
    $data = $fundRows[$i]

    $quarterSheet.Range("A" + $row).Value = $i

    # Columns B-G are stored as text (matching the rest of the workbook,
    # where percentages/amounts are kept as plain strings, not numbers).
    for ($c = 0; $c -lt 6; $c++) {
        $cell = $quarterSheet.Range($headerCols[$c] + $row)
        $cell.NumberFormat = "@"
        $cell.Value = $data[$c]
        $cell.Style = "Normal"
    }

    # Column H (ranking) is a real number.
    $quarterSheet.Range("H" + $row).Value = $data[6]
}

# ---------------------------------------------------------------------
# 3. Rebuild the "总计" sheet: prepend the new 2022-Q1 summary row above
#    the existing 2021-Q4 / 2021-Q3 rows.
# ---------------------------------------------------------------------
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial($xlPasteFormats)

$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 2
$totalSheet.Range("D4").Value = 0.09
$totalSheet.Range("A4").Value = 2

$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 6
$totalSheet.Range("D3").Value = 0.72
$totalSheet.Range("A3").Value = 1

$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 5
$totalSheet.Range("D2").Value = 0.35
$totalSheet.Range("A2").Value = 0
